$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 18 - this pushes the old row 18 (grand total)
# down to row 19, and the old row 19 (footer: timestamp / page / developer)
# down to row 20.
$ws.Rows("18:18").Insert()

# Fill in the new item row (#12 in the list) - a new product sold today.
# A leading apostrophe forces these numeric-looking values to be stored as
# text, matching every other item row in the sheet (counts/prices are text,
# not numbers, even though their column styles use numeric formats).
$ws.Range("A18").Value = 12
$ws.Range("C18").Value = "مناديل سولو سحب"
$ws.Range("H18").Value = "'23:0"
$ws.Range("L18").Value = "'0"
$ws.Range("N18").Value = "'45.00"
$ws.Range("P18").Value = "'45.0000"
$ws.Range("Q18").Value = "'1:0"

# Copy the formatting (styles/number formats) from row 17 (the last item
# row) onto row 18 - applied after the values so the text-forced cells keep
# the normal column style instead of a quote-prefixed default style.
$ws.Range("A17:Q17").Copy()
$ws.Range("A18:Q18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Re-create the merged cell regions for row 18, matching the pattern used
# by every other item row (7-17).
$ws.Range("A18:B18").Merge()
$ws.Range("C18:G18").Merge()
$ws.Range("H18:K18").Merge()
$ws.Range("L18:M18").Merge()
$ws.Range("N18:O18").Merge()

# Row heights: new item row 18, and the grand total row (now pushed to 19).
$ws.Rows("18:18").RowHeight = 24.75
$ws.Rows("19:19").RowHeight = 25.5

# Update the grand total to include the new item's price (675.22 + 45.00).
$ws.Range("P19").Value = 720.22000000000003

# Bump the generated-on timestamp shown in the footer (now row 20).
$ws.Range("A20").Value = "Thursday, 11 September, 2025 11:46 AM"
